# Apply Natmi edit: add M1 sending-cluster block, update values (per Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell contents (keeps header formatting/styles intact)
# and allows shared strings to be rebuilt in the desired order.
$ws.Cells.ClearContents()

# --- Header row (row 1) ---
$ws.Range('A1').Value = 'Sending cluster'
$ws.Range('B1').Value = 'Ligand symbol'
$ws.Range('C1').Value = 'Receptor symbol'
$ws.Range('D1').Value = 'Target cluster'
$ws.Range('E1').Value = 'Ligand-expressing cells'
$ws.Range('F1').Value = 'Ligand detection rate'
$ws.Range('G1').Value = 'Ligand average expression value'
$ws.Range('H1').Value = 'Ligand total expression value'
$ws.Range('I1').Value = 'Ligand derived specificity of average expression value'
$ws.Range('J1').Value = 'Ligand derived specificity of total expression value'
$ws.Range('K1').Value = 'Receptor-expressing cells'
$ws.Range('L1').Value = 'Receptor detection rate'
$ws.Range('M1').Value = 'Receptor average expression value'
$ws.Range('N1').Value = 'Receptor total expression value'
$ws.Range('O1').Value = 'Receptor derived specificity of average expression value'
$ws.Range('P1').Value = 'Receptor derived specificity of total expression value'
$ws.Range('Q1').Value = 'Edge average expression weight'
$ws.Range('R1').Value = 'Edge total expression weight'
$ws.Range('S1').Value = 'Edge average expression derived specificity'
$ws.Range('T1').Value = 'Edge total expression derived specificity'

# --- String columns (A-D), filled column-by-column so the shared-string table
#     is rebuilt with M1 before M2, matching the target workbook ---
# Column A
$ws.Range('A2').Value = 'M1'
$ws.Range('A3').Value = 'M1'
$ws.Range('A4').Value = 'M1'
$ws.Range('A5').Value = 'M1'
$ws.Range('A6').Value = 'M1'
$ws.Range('A7').Value = 'M1'
$ws.Range('A8').Value = 'M2'
$ws.Range('A9').Value = 'M2'
$ws.Range('A10').Value = 'M2'
$ws.Range('A11').Value = 'M2'
$ws.Range('A12').Value = 'M2'
$ws.Range('A13').Value = 'M2'
$ws.Range('A14').Value = 'Neutro'
$ws.Range('A15').Value = 'Neutro'
$ws.Range('A16').Value = 'Neutro'
$ws.Range('A17').Value = 'Neutro'
$ws.Range('A18').Value = 'Neutro'
$ws.Range('A19').Value = 'Neutro'

# Column B
$ws.Range('B2').Value = 'Tnfsf14'
$ws.Range('B3').Value = 'Tnfsf14'
$ws.Range('B4').Value = 'Tnfsf14'
$ws.Range('B5').Value = 'Tnfsf14'
$ws.Range('B6').Value = 'Tnfsf14'
$ws.Range('B7').Value = 'Tnfsf14'
$ws.Range('B8').Value = 'Tnfsf14'
$ws.Range('B9').Value = 'Tnfsf14'
$ws.Range('B10').Value = 'Tnfsf14'
$ws.Range('B11').Value = 'Tnfsf14'
$ws.Range('B12').Value = 'Tnfsf14'
$ws.Range('B13').Value = 'Tnfsf14'
$ws.Range('B14').Value = 'Tnfsf14'
$ws.Range('B15').Value = 'Tnfsf14'
$ws.Range('B16').Value = 'Tnfsf14'
$ws.Range('B17').Value = 'Tnfsf14'
$ws.Range('B18').Value = 'Tnfsf14'
$ws.Range('B19').Value = 'Tnfsf14'

# Column C
$ws.Range('C2').Value = 'Ltbr'
$ws.Range('C3').Value = 'Ltbr'
$ws.Range('C4').Value = 'Ltbr'
$ws.Range('C5').Value = 'Ltbr'
$ws.Range('C6').Value = 'Ltbr'
$ws.Range('C7').Value = 'Ltbr'
$ws.Range('C8').Value = 'Ltbr'
$ws.Range('C9').Value = 'Ltbr'
$ws.Range('C10').Value = 'Ltbr'
$ws.Range('C11').Value = 'Ltbr'
$ws.Range('C12').Value = 'Ltbr'
$ws.Range('C13').Value = 'Ltbr'
$ws.Range('C14').Value = 'Ltbr'
$ws.Range('C15').Value = 'Ltbr'
$ws.Range('C16').Value = 'Ltbr'
$ws.Range('C17').Value = 'Ltbr'
$ws.Range('C18').Value = 'Ltbr'
$ws.Range('C19').Value = 'Ltbr'

# Column D
$ws.Range('D2').Value = 'ECs'
$ws.Range('D3').Value = 'FAPs'
$ws.Range('D4').Value = 'M1'
$ws.Range('D5').Value = 'M2'
$ws.Range('D6').Value = 'Neutro'
$ws.Range('D7').Value = 'sCs'
$ws.Range('D8').Value = 'ECs'
$ws.Range('D9').Value = 'FAPs'
$ws.Range('D10').Value = 'M1'
$ws.Range('D11').Value = 'M2'
$ws.Range('D12').Value = 'Neutro'
$ws.Range('D13').Value = 'sCs'
$ws.Range('D14').Value = 'ECs'
$ws.Range('D15').Value = 'FAPs'
$ws.Range('D16').Value = 'M1'
$ws.Range('D17').Value = 'M2'
$ws.Range('D18').Value = 'Neutro'
$ws.Range('D19').Value = 'sCs'

# --- Numeric columns (E-T) ---
# Row 2
$ws.Range('E2').Value2 = 2
$ws.Range('F2').Value2 = 0.6666666666666666
$ws.Range('G2').Value2 = 2.786831
$ws.Range('H2').Value2 = 8.360493
$ws.Range('I2').Value2 = 0.1227769702371957
$ws.Range('J2').Value2 = 0.1227769702371957
$ws.Range('K2').Value2 = 3
$ws.Range('L2').Value2 = 1
$ws.Range('M2').Value2 = 8.226432
$ws.Range('N2').Value2 = 24.679296
$ws.Range('O2').Value2 = 0.1046851189010213
$ws.Range('P2').Value2 = 0.1046851189010213
$ws.Range('Q2').Value2 = 22.925675716992
$ws.Range('R2').Value2 = 206.331081452928
$ws.Range('S2').Value2 = 0.01285292172758798
$ws.Range('T2').Value2 = 0.01285292172758798

# Row 3
$ws.Range('E3').Value2 = 2
$ws.Range('F3').Value2 = 0.6666666666666666
$ws.Range('G3').Value2 = 2.786831
$ws.Range('H3').Value2 = 8.360493
$ws.Range('I3').Value2 = 0.1227769702371957
$ws.Range('J3').Value2 = 0.1227769702371957
$ws.Range('K3').Value2 = 3
$ws.Range('L3').Value2 = 1
$ws.Range('M3').Value2 = 20.15320433333333
$ws.Range('N3').Value2 = 60.45961299999999
$ws.Range('O3').Value2 = 0.2564587650966515
$ws.Range('P3').Value2 = 0.2564587650966514
$ws.Range('Q3').Value2 = 56.16357458546766
$ws.Range('R3').Value2 = 505.4721712692089
$ws.Range('S3').Value2 = 0.03148723016933954
$ws.Range('T3').Value2 = 0.03148723016933953

# Row 4
$ws.Range('E4').Value2 = 2
$ws.Range('F4').Value2 = 0.6666666666666666
$ws.Range('G4').Value2 = 2.786831
$ws.Range('H4').Value2 = 8.360493
$ws.Range('I4').Value2 = 0.1227769702371957
$ws.Range('J4').Value2 = 0.1227769702371957
$ws.Range('K4').Value2 = 3
$ws.Range('L4').Value2 = 1
$ws.Range('M4').Value2 = 15.64009466666667
$ws.Range('N4').Value2 = 46.920284
$ws.Range('O4').Value2 = 0.1990273753922999
$ws.Range('P4').Value2 = 0.1990273753922999
$ws.Range('Q4').Value2 = 43.58630066000133
$ws.Range('R4').Value2 = 392.276705940012
$ws.Range('S4').Value2 = 0.02443597814492758
$ws.Range('T4').Value2 = 0.02443597814492758

# Row 5
$ws.Range('E5').Value2 = 2
$ws.Range('F5').Value2 = 0.6666666666666666
$ws.Range('G5').Value2 = 2.786831
$ws.Range('H5').Value2 = 8.360493
$ws.Range('I5').Value2 = 0.1227769702371957
$ws.Range('J5').Value2 = 0.1227769702371957
$ws.Range('K5').Value2 = 3
$ws.Range('L5').Value2 = 1
$ws.Range('M5').Value2 = 13.81253266666666
$ws.Range('N5').Value2 = 41.43759799999999
$ws.Range('O5').Value2 = 0.1757708110313487
$ws.Range('P5').Value2 = 0.1757708110313487
$ws.Range('Q5').Value2 = 38.49319422397932
$ws.Range('R5').Value2 = 346.4387480158139
$ws.Range('S5').Value2 = 0.02158060763456365
$ws.Range('T5').Value2 = 0.02158060763456365

# Row 6
$ws.Range('E6').Value2 = 2
$ws.Range('F6').Value2 = 0.6666666666666666
$ws.Range('G6').Value2 = 2.786831
$ws.Range('H6').Value2 = 8.360493
$ws.Range('I6').Value2 = 0.1227769702371957
$ws.Range('J6').Value2 = 0.1227769702371957
$ws.Range('K6').Value2 = 3
$ws.Range('L6').Value2 = 1
$ws.Range('M6').Value2 = 13.92712666666667
$ws.Range('N6').Value2 = 41.78138
$ws.Range('O6').Value2 = 0.1772290722210533
$ws.Range('P6').Value2 = 0.1772290722210533
$ws.Range('Q6').Value2 = 38.81254833559333
$ws.Range('R6').Value2 = 349.31293502034
$ws.Range('S6').Value2 = 0.02175964852525007
$ws.Range('T6').Value2 = 0.02175964852525007

# Row 7
$ws.Range('E7').Value2 = 2
$ws.Range('F7').Value2 = 0.6666666666666666
$ws.Range('G7').Value2 = 2.786831
$ws.Range('H7').Value2 = 8.360493
$ws.Range('I7').Value2 = 0.1227769702371957
$ws.Range('J7').Value2 = 0.1227769702371957
$ws.Range('K7').Value2 = 3
$ws.Range('L7').Value2 = 1
$ws.Range('M7').Value2 = 6.823240000000001
$ws.Range('N7').Value2 = 20.46972
$ws.Range('O7').Value2 = 0.08682885735762533
$ws.Range('P7').Value2 = 0.08682885735762533
$ws.Range('Q7').Value2 = 19.01521675244
$ws.Range('R7').Value2 = 171.13695077196
$ws.Range('S7').Value2 = 0.01066058403552687
$ws.Range('T7').Value2 = 0.01066058403552687

# Row 8
$ws.Range('E8').Value2 = 3
$ws.Range('F8').Value2 = 1
$ws.Range('G8').Value2 = 2.556772666666667
$ws.Range('H8').Value2 = 7.670318
$ws.Range('I8').Value2 = 0.1126414919306584
$ws.Range('J8').Value2 = 0.1126414919306584
$ws.Range('K8').Value2 = 3
$ws.Range('L8').Value2 = 1
$ws.Range('M8').Value2 = 8.226432
$ws.Range('N8').Value2 = 24.679296
$ws.Range('O8').Value2 = 0.1046851189010213
$ws.Range('P8').Value2 = 0.1046851189010213
$ws.Range('Q8').Value2 = 21.033116481792
$ws.Range('R8').Value2 = 189.298048336128
$ws.Range('S8').Value2 = 0.01179188797594941
$ws.Range('T8').Value2 = 0.01179188797594941

# Row 9
$ws.Range('E9').Value2 = 3
$ws.Range('F9').Value2 = 1
$ws.Range('G9').Value2 = 2.556772666666667
$ws.Range('H9').Value2 = 7.670318
$ws.Range('I9').Value2 = 0.1126414919306584
$ws.Range('J9').Value2 = 0.1126414919306584
$ws.Range('K9').Value2 = 3
$ws.Range('L9').Value2 = 1
$ws.Range('M9').Value2 = 20.15320433333333
$ws.Range('N9').Value2 = 60.45961299999999
$ws.Range('O9').Value2 = 0.2564587650966515
$ws.Range('P9').Value2 = 0.2564587650966514
$ws.Range('Q9').Value2 = 51.52716198521489
$ws.Range('R9').Value2 = 463.7444578669339
$ws.Range('S9').Value2 = 0.0288878979191811
$ws.Range('T9').Value2 = 0.02888789791918109

# Row 10
$ws.Range('E10').Value2 = 3
$ws.Range('F10').Value2 = 1
$ws.Range('G10').Value2 = 2.556772666666667
$ws.Range('H10').Value2 = 7.670318
$ws.Range('I10').Value2 = 0.1126414919306584
$ws.Range('J10').Value2 = 0.1126414919306584
$ws.Range('K10').Value2 = 3
$ws.Range('L10').Value2 = 1
$ws.Range('M10').Value2 = 15.64009466666667
$ws.Range('N10').Value2 = 46.920284
$ws.Range('O10').Value2 = 0.1990273753922999
$ws.Range('P10').Value2 = 0.1990273753922999
$ws.Range('Q10').Value2 = 39.98816654781245
$ws.Range('R10').Value2 = 359.893498930312
$ws.Range('S10').Value2 = 0.02241874049923188
$ws.Range('T10').Value2 = 0.02241874049923188

# Row 11
$ws.Range('E11').Value2 = 3
$ws.Range('F11').Value2 = 1
$ws.Range('G11').Value2 = 2.556772666666667
$ws.Range('H11').Value2 = 7.670318
$ws.Range('I11').Value2 = 0.1126414919306584
$ws.Range('J11').Value2 = 0.1126414919306584
$ws.Range('K11').Value2 = 3
$ws.Range('L11').Value2 = 1
$ws.Range('M11').Value2 = 13.81253266666666
$ws.Range('N11').Value2 = 41.43759799999999
$ws.Range('O11').Value2 = 0.1757708110313487
$ws.Range('P11').Value2 = 0.1757708110313487
$ws.Range('Q11').Value2 = 35.31550597957377
$ws.Range('R11').Value2 = 317.839553816164
$ws.Range('S11').Value2 = 0.01979908639243296
$ws.Range('T11').Value2 = 0.01979908639243295

# Row 12
$ws.Range('E12').Value2 = 3
$ws.Range('F12').Value2 = 1
$ws.Range('G12').Value2 = 2.556772666666667
$ws.Range('H12').Value2 = 7.670318
$ws.Range('I12').Value2 = 0.1126414919306584
$ws.Range('J12').Value2 = 0.1126414919306584
$ws.Range('K12').Value2 = 3
$ws.Range('L12').Value2 = 1
$ws.Range('M12').Value2 = 13.92712666666667
$ws.Range('N12').Value2 = 41.78138
$ws.Range('O12').Value2 = 0.1772290722210533
$ws.Range('P12').Value2 = 0.1772290722210533
$ws.Range('Q12').Value2 = 35.60849678653778
$ws.Range('R12').Value2 = 320.47647107884
$ws.Range('S12').Value2 = 0.01996334710846586
$ws.Range('T12').Value2 = 0.01996334710846585

# Row 13
$ws.Range('E13').Value2 = 3
$ws.Range('F13').Value2 = 1
$ws.Range('G13').Value2 = 2.556772666666667
$ws.Range('H13').Value2 = 7.670318
$ws.Range('I13').Value2 = 0.1126414919306584
$ws.Range('J13').Value2 = 0.1126414919306584
$ws.Range('K13').Value2 = 3
$ws.Range('L13').Value2 = 1
$ws.Range('M13').Value2 = 6.823240000000001
$ws.Range('N13').Value2 = 20.46972
$ws.Range('O13').Value2 = 0.08682885735762533
$ws.Range('P13').Value2 = 0.08682885735762533
$ws.Range('Q13').Value2 = 17.44547353010667
$ws.Range('R13').Value2 = 157.00926177096
$ws.Range('S13').Value2 = 0.009780532035397246
$ws.Range('T13').Value2 = 0.009780532035397244

# Row 14
$ws.Range('E14').Value2 = 3
$ws.Range('F14').Value2 = 1
$ws.Range('G14').Value2 = 17.35471666666666
$ws.Range('H14').Value2 = 52.06415
$ws.Range('I14').Value2 = 0.7645815378321459
$ws.Range('J14').Value2 = 0.7645815378321459
$ws.Range('K14').Value2 = 3
$ws.Range('L14').Value2 = 1
$ws.Range('M14').Value2 = 8.226432
$ws.Range('N14').Value2 = 24.679296
$ws.Range('O14').Value2 = 0.1046851189010213
$ws.Range('P14').Value2 = 0.1046851189010213
$ws.Range('Q14').Value2 = 142.7673965376
$ws.Range('R14').Value2 = 1284.9065688384
$ws.Range('S14').Value2 = 0.0800403091974839
$ws.Range('T14').Value2 = 0.0800403091974839

# Row 15
$ws.Range('E15').Value2 = 3
$ws.Range('F15').Value2 = 1
$ws.Range('G15').Value2 = 17.35471666666666
$ws.Range('H15').Value2 = 52.06415
$ws.Range('I15').Value2 = 0.7645815378321459
$ws.Range('J15').Value2 = 0.7645815378321459
$ws.Range('K15').Value2 = 3
$ws.Range('L15').Value2 = 1
$ws.Range('M15').Value2 = 20.15320433333333
$ws.Range('N15').Value2 = 60.45961299999999
$ws.Range('O15').Value2 = 0.2564587650966515
$ws.Range('P15').Value2 = 0.2564587650966514
$ws.Range('Q15').Value2 = 349.7531511304388
$ws.Range('R15').Value2 = 3147.778360173949
$ws.Range('S15').Value2 = 0.1960836370081309
$ws.Range('T15').Value2 = 0.1960836370081308

# Row 16
$ws.Range('E16').Value2 = 3
$ws.Range('F16').Value2 = 1
$ws.Range('G16').Value2 = 17.35471666666666
$ws.Range('H16').Value2 = 52.06415
$ws.Range('I16').Value2 = 0.7645815378321459
$ws.Range('J16').Value2 = 0.7645815378321459
$ws.Range('K16').Value2 = 3
$ws.Range('L16').Value2 = 1
$ws.Range('M16').Value2 = 15.64009466666667
$ws.Range('N16').Value2 = 46.920284
$ws.Range('O16').Value2 = 0.1990273753922999
$ws.Range('P16').Value2 = 0.1990273753922999
$ws.Range('Q16').Value2 = 271.4294115798444
$ws.Range('R16').Value2 = 2442.8647042186
$ws.Range('S16').Value2 = 0.1521726567481405
$ws.Range('T16').Value2 = 0.1521726567481405

# Row 17
$ws.Range('E17').Value2 = 3
$ws.Range('F17').Value2 = 1
$ws.Range('G17').Value2 = 17.35471666666666
$ws.Range('H17').Value2 = 52.06415
$ws.Range('I17').Value2 = 0.7645815378321459
$ws.Range('J17').Value2 = 0.7645815378321459
$ws.Range('K17').Value2 = 3
$ws.Range('L17').Value2 = 1
$ws.Range('M17').Value2 = 13.81253266666666
$ws.Range('N17').Value2 = 41.43759799999999
$ws.Range('O17').Value2 = 0.1757708110313487
$ws.Range('P17').Value2 = 0.1757708110313487
$ws.Range('Q17').Value2 = 239.7125908790777
$ws.Range('R17').Value2 = 2157.4133179117
$ws.Range('S17').Value2 = 0.1343911170043521
$ws.Range('T17').Value2 = 0.1343911170043521

# Row 18
$ws.Range('E18').Value2 = 3
$ws.Range('F18').Value2 = 1
$ws.Range('G18').Value2 = 17.35471666666666
$ws.Range('H18').Value2 = 52.06415
$ws.Range('I18').Value2 = 0.7645815378321459
$ws.Range('J18').Value2 = 0.7645815378321459
$ws.Range('K18').Value2 = 3
$ws.Range('L18').Value2 = 1
$ws.Range('M18').Value2 = 13.92712666666667
$ws.Range('N18').Value2 = 41.78138
$ws.Range('O18').Value2 = 0.1772290722210533
$ws.Range('P18').Value2 = 0.1772290722210533
$ws.Range('Q18').Value2 = 241.7013372807778
$ws.Range('R18').Value2 = 2175.312035527
$ws.Range('S18').Value2 = 0.1355060765873374
$ws.Range('T18').Value2 = 0.1355060765873374

# Row 19
$ws.Range('E19').Value2 = 3
$ws.Range('F19').Value2 = 1
$ws.Range('G19').Value2 = 17.35471666666666
$ws.Range('H19').Value2 = 52.06415
$ws.Range('I19').Value2 = 0.7645815378321459
$ws.Range('J19').Value2 = 0.7645815378321459
$ws.Range('K19').Value2 = 3
$ws.Range('L19').Value2 = 1
$ws.Range('M19').Value2 = 6.823240000000001
$ws.Range('N19').Value2 = 20.46972
$ws.Range('O19').Value2 = 0.08682885735762533
$ws.Range('P19').Value2 = 0.08682885735762533
$ws.Range('Q19').Value2 = 118.4153969486667
$ws.Range('R19').Value2 = 1065.738572538
$ws.Range('S19').Value2 = 0.0663877412867012
$ws.Range('T19').Value2 = 0.0663877412867012
